# Natmi following Dr Hou advice
# Re-derive the LR-pair table (Adm -> Calcrl) including "ECs" as a sending
# cluster (it previously only appeared as a target cluster), giving the full
# 3x3 cross product of {ECs, FAPs, sCs} sending/target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Adm", "Calcrl", "ECs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 3, 1, 46.06383033333333, 138.191491, 0.6011809800814781, 0.601180980081478, 497.2048313200308, 4474.843481880278, 0.1790282524409286, 0.1790282524409285),
  @("ECs", "Adm", "Calcrl", "FAPs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 3, 1, 24.118612, 72.35583600000001, 0.3147730159528759, 0.3147730159528759, 260.332028933676, 2342.988260403084, 0.0937376012028297, 0.0937376012028297),
  @("ECs", "Adm", "Calcrl", "sCs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 3, 1, 6.439792666666666, 19.319378, 0.0840460039656461, 0.0840460039656461, 69.50998220069799, 625.5898398062819, 0.025028418584656, 0.025028418584656),
  @("FAPs", "Adm", "Calcrl", "ECs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 3, 1, 46.06383033333333, 138.191491, 0.6011809800814781, 0.601180980081478, 1096.509067751188, 9868.581609760688, 0.3948193778888936, 0.3948193778888935),
  @("FAPs", "Adm", "Calcrl", "FAPs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 3, 1, 24.118612, 72.35583600000001, 0.3147730159528759, 0.3147730159528759, 574.1223985977388, 5167.101587379649, 0.2067239158462428, 0.2067239158462428),
  @("FAPs", "Adm", "Calcrl", "sCs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 3, 1, 6.439792666666666, 19.319378, 0.0840460039656461, 0.0840460039656461, 153.2936145852338, 1379.642531267104, 0.05519634203208921, 0.05519634203208921),
  @("sCs", "Adm", "Calcrl", "ECs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 3, 1, 46.06383033333333, 138.191491, 0.6011809800814781, 0.601180980081478, 75.9113344815099, 683.202010333589, 0.02733334975165598, 0.02733334975165598),
  @("sCs", "Adm", "Calcrl", "FAPs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 3, 1, 24.118612, 72.35583600000001, 0.3147730159528759, 0.3147730159528759, 39.74649979198268, 357.7184981278441, 0.01431149890380343, 0.01431149890380343),
  @("sCs", "Adm", "Calcrl", "sCs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 3, 1, 6.439792666666666, 19.319378, 0.0840460039656461, 0.0840460039656461, 10.61251857636245, 95.51266718726201, 0.003821243348900896, 0.003821243348900897)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

Write-Output "rows written"
